$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 91 ---
$ws.Range("A90").Copy()
$ws.Range("A91").PasteSpecial(-4122)  # xlPasteFormats - reuse the date/time style
$ws.Range("A91").Value = 45454.2916666667
$ws.Range("B91").Value = 4399
$ws.Range("C91").Value = 0.720000028610229
$ws.Range("D91").Value = 0.714999973773956
$ws.Range("E91").Value = 0.720000028610229
$ws.Range("F91").Value = 0.714999973773956

$g91 = $ws.Cells.Item(91, 7)
$g91.NumberFormat = "@"
$g91.Value = "0.714999973773956"
$g91.Style = "Normal"

$ws.Range("H91").Value = "BWZ.MI"

# --- Row 92 ---
$ws.Range("A90").Copy()
$ws.Range("A92").PasteSpecial(-4122)  # xlPasteFormats - reuse the date/time style
$ws.Range("A92").Value = 45455.6412962963
$ws.Range("B92").Value = 16624
$ws.Range("C92").Value = 0.730000019073486
$ws.Range("D92").Value = 0.709999978542328
$ws.Range("E92").Value = 0.709999978542328
$ws.Range("F92").Value = 0.720000028610229

$g92 = $ws.Cells.Item(92, 7)
$g92.NumberFormat = "@"
$g92.Value = "0.720000028610229"
$g92.Style = "Normal"

$ws.Range("H92").Value = "BWZ.MI"
